$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift existing rows 16-30 down to 17-31 (bottom-up so we never clobber
#    a row before it has been copied).
# ---------------------------------------------------------------------------
for ($r = 30; $r -ge 16; $r--) {
    $src = $ws.Range("A" + $r + ":P" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":P" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# Row 31 is brand new territory (previously empty) - copy formatting from the
# row above it so column A keeps the bold/bordered "index" style.
$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Populate the newly inserted row 16 with its own fresh data.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value2 = "ibes_qoq_2|fwdepsqcut|q_2｜ibes_qoq_tune10_ind"
$ws.Range("B16").Value2 = 0.003103990424644728
$ws.Range("C16").Value2 = 0.00254141838673398
$ws.Range("D16").Value2 = 0.00002156002471405648
$ws.Range("E16").Value2 = 0.00002106152912726658
$ws.Range("F16").Value2 = -0.02262657763782516
$ws.Range("G16").Value2 = 0.001017868166272584
$ws.Range("H16").Value2 = 0.212555477012045
$ws.Range("I16").Value2 = 11542
$ws.Range("J16").Value2 = 0.000163888349337522
$ws.Range("K16").Value2 = 0.001597174444307401
$ws.Range("L16").Value2 = 0.004589584286176587
$ws.Range("M16").Value2 = 0.00510134746390362
$ws.Range("N16").Value2 = 0.001583264354559398
$ws.Range("O16").Value2 = -0.005012774390606101
$ws.Range("P16").Value2 = 0.0002091701444985722

# ---------------------------------------------------------------------------
# 3. Re-derive the "_org" hyperparameter-space columns (J:P) for every
#    pre-existing data row (2-15, and the shifted-down 17-31 - row 16 is the
#    freshly inserted row and already carries its final values from step 2).
#    The column headers for J:P were re-ordered, which amounts to this
#    permutation of each row's existing J:P values:
#       new J = old P      new K = old O      new L = old K
#       new M = old L      new N = old J      new O = old M
#       new P = old N
# ---------------------------------------------------------------------------
$dataRows = @(2..15) + @(17..31)
foreach ($r in $dataRows) {
    $row = $ws.Range("J" + $r + ":P" + $r).Value2

    $oldJ = $row[1,1]
    $oldK = $row[1,2]
    $oldL = $row[1,3]
    $oldM = $row[1,4]
    $oldN = $row[1,5]
    $oldO = $row[1,6]
    $oldP = $row[1,7]

    $ws.Range("J" + $r).Value2 = $oldP
    $ws.Range("K" + $r).Value2 = $oldO
    $ws.Range("L" + $r).Value2 = $oldK
    $ws.Range("M" + $r).Value2 = $oldL
    $ws.Range("N" + $r).Value2 = $oldJ
    $ws.Range("O" + $r).Value2 = $oldM
    $ws.Range("P" + $r).Value2 = $oldN
}

# ---------------------------------------------------------------------------
# 4. Header row labels (J1:P1) swap meaning the same way the shared strings
#    were re-ordered, so make sure the visible header text matches.
# ---------------------------------------------------------------------------
$ws.Range("J1").Value2 = "consensus_mse_org"
$ws.Range("K1").Value2 = "consensus_medae_org"
$ws.Range("L1").Value2 = "consensus_mae_org"
$ws.Range("M1").Value2 = "lgbm_mae_org"
$ws.Range("N1").Value2 = "lgbm_medae_org"
$ws.Range("O1").Value2 = "lgbm_r2_org"
$ws.Range("P1").Value2 = "lgbm_mse_org"
